# Append the latest run-log row (row 5) to the Nalco run log sheet,
# mirroring the formatting of the previous "SKIPPED" row (row 4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of row 4 onto the new row 5 first, so the freshly
# written cells inherit the same style (centered, default font/fill/border)
# as the existing log rows instead of Excel's generic default style.
$ws.Range("A4:H4").Copy()
$ws.Range("A5:H5").PasteSpecial(-4122) # xlPasteFormats

$ws.Cells.Item(5, 1).Value = "2025-11-07 06:45:07 UTC"
$ws.Cells.Item(5, 2).Value = "2025-11-07 12:15:07 IST"
$ws.Cells.Item(5, 3).Value = "SKIPPED"
$ws.Cells.Item(5, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item(5, 5).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-11-2025.pdf"
$ws.Cells.Item(5, 6).Value = ""
$ws.Cells.Item(5, 7).Value = 0
$ws.Cells.Item(5, 8).Value = ""

$wb.Save()
